# Generate Report for Handoff
# Update the "Latest Handoff Datetime" for the e713a193-... file row (row 5)
# on both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-03-09 05:27:47"

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-03-09 05:27:56"
